# Apply cryptos list refresh (GitHub Actions daily update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains numbers formatted with literal "." thousands
# separators, stored as plain text in the sheet (t="inlineStr"). Several of the
# refreshed values now look like ordinary decimals (e.g. "325.13") which Excel
# would otherwise auto-coerce to a Number on assignment. Force the whole Price
# column to Text for the duration of the writes, then restore the default
# (unstyled) formatting so the saved cells match the original look & feel.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "29.436.80"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.908.61"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "325.13"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "23.41"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").Value = "1.914.21"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "6.012"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "7.157"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "90.26"
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.06791"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "17.66"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "29.457.32"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "5.626"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "11.71"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").Value = "2.175"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "2.158.40"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "156.30"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "6.489"
$ws.Range("E27").Value = "  +8.30%  "
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "120.35"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "1.027"
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").Value = "0.09517"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "5.508"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "3.562"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "1.389"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").Value = "0.02268"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "0.06101"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "1.174"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "10.81"
$ws.Range("E39").Value = "  +7.04%  "
$ws.Range("D40").Value = "0.5945"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "7.972"
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").Value = "0.1854"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "1.275"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "2.376"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").Value = "12.55"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("D46").Value = "0.07602"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").Value = "0.5561"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").Value = "116.50"
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").Value = "72.53"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "2.404"
$ws.Range("E51").Value = "  +2.41%  "

$priceCol.NumberFormat = "General"
$priceCol.Style = "Normal"

